# Apply corrected IFRS financial figures (단위 변경: 백만원 -> 십억원 등 보정)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3146
$ws.Range("E2").Value = 97
$ws.Range("F2").Value = 97
$ws.Range("G2").Value = 72
$ws.Range("H2").Value = 41
$ws.Range("I2").Value = 41
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2996
$ws.Range("L2").Value = 1335
$ws.Range("M2").Value = 1661
$ws.Range("N2").Value = 1641
$ws.Range("O2").Value = 20
$ws.Range("P2").Value = 377
$ws.Range("Q2").Value = 123
$ws.Range("R2").Value = -198
$ws.Range("S2").Value = -25
$ws.Range("T2").Value = 48
$ws.Range("U2").Value = 74
$ws.Range("V2").Value = 532
$ws.Range("W2").Value = 3.09
$ws.Range("X2").Value = 1.29
$ws.Range("Y2").Value = 2.49
$ws.Range("Z2").Value = 1.36
$ws.Range("AA2").Value = 80.37
$ws.Range("AB2").Value = 323.51
$ws.Range("AC2").Value = 539
$ws.Range("AD2").Value = 39.42
$ws.Range("AE2").Value = 21782
$ws.Range("AF2").Value = 0.98
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 0.47
$ws.Range("AI2").Value = 18.55
$ws.Range("AJ2").Value = 7533015

# Row 3
$ws.Range("D3").Value = 2909
$ws.Range("E3").Value = 71
$ws.Range("F3").Value = 71
$ws.Range("G3").Value = -1
$ws.Range("H3").Value = -34
$ws.Range("I3").Value = -32
$ws.Range("J3").Value = -3
$ws.Range("K3").Value = 2965
$ws.Range("L3").Value = 1340
$ws.Range("M3").Value = 1625
$ws.Range("N3").Value = 1607
$ws.Range("O3").Value = 18
$ws.Range("P3").Value = 377
$ws.Range("Q3").Value = 131
$ws.Range("R3").Value = -70
$ws.Range("S3").Value = -24
$ws.Range("T3").Value = 73
$ws.Range("U3").Value = 57
$ws.Range("V3").Value = 532
$ws.Range("W3").Value = 2.46
$ws.Range("X3").Value = -1.17
$ws.Range("Y3").Value = -1.94
$ws.Range("Z3").Value = -1.14
$ws.Range("AA3").Value = 82.45
$ws.Range("AB3").Value = 310.58
$ws.Range("AC3").Value = -418
$ws.Range("AD3").Value = -65.95999999999999
$ws.Range("AE3").Value = 21339
$ws.Range("AF3").Value = 1.29
$ws.Range("AG3").Value = 50
$ws.Range("AH3").Value = 0.18
$ws.Range("AI3").Value = -11.95
$ws.Range("AJ3").Value = 7533015

# Row 4
$ws.Range("D4").Value = 3593
$ws.Range("E4").Value = 253
$ws.Range("F4").Value = 253
$ws.Range("G4").Value = 247
$ws.Range("H4").Value = 188
$ws.Range("I4").Value = 183
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 3267
$ws.Range("L4").Value = 1492
$ws.Range("M4").Value = 1775
$ws.Range("N4").Value = 1753
$ws.Range("O4").Value = 23
$ws.Range("P4").Value = 377
$ws.Range("Q4").Value = 372
$ws.Range("R4").Value = -177
$ws.Range("S4").Value = -29
$ws.Range("T4").Value = 192
$ws.Range("U4").Value = 180
$ws.Range("V4").Value = 520
$ws.Range("W4").Value = 7.03
$ws.Range("X4").Value = 5.23
$ws.Range("Y4").Value = 10.9
$ws.Range("Z4").Value = 6.03
$ws.Range("AA4").Value = 84.01000000000001
$ws.Range("AB4").Value = 355.43
$ws.Range("AC4").Value = 2431
$ws.Range("AD4").Value = 16.68
$ws.Range("AE4").Value = 23270
$ws.Range("AF4").Value = 1.74
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 0.37
$ws.Range("AI4").Value = 6.17
$ws.Range("AJ4").Value = 7533015

# Row 5
$ws.Range("D5").Value = 4585
$ws.Range("E5").Value = 433
$ws.Range("F5").Value = 433
$ws.Range("G5").Value = 357
$ws.Range("H5").Value = 291
$ws.Range("I5").Value = 287
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 3982
$ws.Range("L5").Value = 1951
$ws.Range("M5").Value = 2031
$ws.Range("N5").Value = 2001
$ws.Range("O5").Value = 30
$ws.Range("P5").Value = 377
$ws.Range("Q5").Value = 436
$ws.Range("R5").Value = -456
$ws.Range("S5").Value = 121
$ws.Range("T5").Value = 456
$ws.Range("U5").Value = -20
$ws.Range("V5").Value = 666
$ws.Range("W5").Value = 9.449999999999999
$ws.Range("X5").Value = 6.35
$ws.Range("Y5").Value = 15.28
$ws.Range("Z5").Value = 8.029999999999999
$ws.Range("AA5").Value = 96.04000000000001
$ws.Range("AB5").Value = 421.61
$ws.Range("AC5").Value = 3807
$ws.Range("AD5").Value = 25.42
$ws.Range("AE5").Value = 26562
$ws.Range("AF5").Value = 3.64
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 0.26
$ws.Range("AI5").Value = 6.57
$ws.Range("AJ5").Value = 7533015

# Row 6
$ws.Range("D6").Value = 4694
$ws.Range("E6").Value = 552
$ws.Range("F6").Value = 552
$ws.Range("G6").Value = 539
$ws.Range("H6").Value = 353
$ws.Range("I6").Value = 348
$ws.Range("K6").Value = 4121
$ws.Range("L6").Value = 1777
$ws.Range("M6").Value = 2344
$ws.Range("N6").Value = 2309
$ws.Range("P6").Value = 377
$ws.Range("Q6").Value = 380
$ws.Range("R6").Value = -230
$ws.Range("S6").Value = -112
$ws.Range("T6").Value = 356
$ws.Range("U6").Value = 24
$ws.Range("V6").Value = 593
$ws.Range("W6").Value = 11.76
$ws.Range("X6").Value = 7.52
$ws.Range("Y6").Value = 16.16
$ws.Range("Z6").Value = 8.710000000000001
$ws.Range("AA6").Value = 75.81
$ws.Range("AB6").Value = 504.19
$ws.Range("AC6").Value = 4623
$ws.Range("AD6").Value = 11.25
$ws.Range("AE6").Value = 30654
$ws.Range("AF6").Value = 1.7
$ws.Range("AG6").Value = 400
$ws.Range("AH6").Value = 0.77
$ws.Range("AI6").Value = 8.65
$ws.Range("AJ6").Value = 7533015

# Row 7
$ws.Range("D7").Value = 5414
$ws.Range("E7").Value = 811
$ws.Range("G7").Value = 817
$ws.Range("H7").Value = 629
$ws.Range("I7").Value = 622
$ws.Range("K7").Value = 4998
$ws.Range("L7").Value = 2073
$ws.Range("M7").Value = 2925
$ws.Range("N7").Value = 2895
$ws.Range("P7").Value = 378
$ws.Range("Q7").Value = 631
$ws.Range("R7").Value = -198
$ws.Range("S7").Value = -85
$ws.Range("T7").Value = 136
$ws.Range("U7").Value = 322
$ws.Range("W7").Value = 14.97
$ws.Range("X7").Value = 11.62
$ws.Range("Y7").Value = 23.92
$ws.Range("Z7").Value = 13.8
$ws.Range("AA7").Value = 70.89
$ws.Range("AC7").Value = 8262
$ws.Range("AD7").Value = 12.28
$ws.Range("AE7").Value = 38433
$ws.Range("AF7").Value = 2.64
$ws.Range("AG7").Value = 410
$ws.Range("AH7").Value = 0.4
$ws.Range("AI7").Value = 4.96

# Row 8
$ws.Range("D8").Value = 5948
$ws.Range("E8").Value = 947
$ws.Range("G8").Value = 960
$ws.Range("H8").Value = 726
$ws.Range("I8").Value = 716
$ws.Range("K8").Value = 5474
$ws.Range("L8").Value = 1863
$ws.Range("M8").Value = 3611
$ws.Range("N8").Value = 3594
$ws.Range("P8").Value = 378
$ws.Range("Q8").Value = 789
$ws.Range("R8").Value = -296
$ws.Range("S8").Value = -76
$ws.Range("T8").Value = 206
$ws.Range("U8").Value = 511
$ws.Range("W8").Value = 15.93
$ws.Range("X8").Value = 12.2
$ws.Range("Y8").Value = 22.08
$ws.Range("Z8").Value = 13.86
$ws.Range("AA8").Value = 51.58
$ws.Range("AC8").Value = 9510
$ws.Range("AD8").Value = 10.67
$ws.Range("AE8").Value = 47713
$ws.Range("AF8").Value = 2.13
$ws.Range("AG8").Value = 420
$ws.Range("AH8").Value = 0.41
$ws.Range("AI8").Value = 4.42

# Row 9
$ws.Range("D9").Value = 6335
$ws.Range("E9").Value = 1004
$ws.Range("G9").Value = 1024
$ws.Range("H9").Value = 769
$ws.Range("I9").Value = 760
$ws.Range("K9").Value = 6160
$ws.Range("L9").Value = 1842
$ws.Range("M9").Value = 4317
$ws.Range("N9").Value = 4324
$ws.Range("P9").Value = 378
$ws.Range("Q9").Value = 871
$ws.Range("R9").Value = -282
$ws.Range("S9").Value = -76
$ws.Range("T9").Value = 158
$ws.Range("U9").Value = 620
$ws.Range("W9").Value = 15.84
$ws.Range("X9").Value = 12.14
$ws.Range("Y9").Value = 19.18
$ws.Range("Z9").Value = 13.23
$ws.Range("AA9").Value = 42.68
$ws.Range("AC9").Value = 10082
$ws.Range("AD9").Value = 10.07
$ws.Range("AE9").Value = 57399
$ws.Range("AF9").Value = 1.77
$ws.Range("AG9").Value = 425
$ws.Range("AH9").Value = 0.42
$ws.Range("AI9").Value = 4.21

